# Add explanation rows (row 2) to "CbCR Notifications" and "TP Forms" sheets,
# pushing the existing example rows down by one, and re-pointing the dropdown
# data validations so they start at row 3 instead of row 2.

$wb = $excel.ActiveWorkbook

# A sheet that already has the correct "explanation row" formatting (yellow
# fill, italic 9pt font, top-aligned wrapped text) on row 2 - use it as the
# formatting source so we reuse the existing style instead of inventing a
# new one.
$formatSource = $wb.Worksheets.Item("MF Requirements")

function Add-ExplanationRow {
    param(
        [string]$SheetName,
        [string]$LastCol,
        [string[]]$Texts,
        [hashtable[]]$Validations
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Insert a new blank row above the current row 2 - this pushes all the
    # existing example rows down by one automatically.
    $ws.Rows.Item(2).Insert()

    # Copy the explanation-row formatting (fill/font/alignment) from the
    # known-good sheet so the new row matches rows already in the workbook.
    $formatSource.Range("A2:" + $LastCol + "2").Copy()
    $ws.Range("A2:" + $LastCol + "2").PasteSpecial(-4122)

    # Explanation rows use a taller, wrapped row height.
    $ws.Rows.Item(2).RowHeight = 60

    # Fill in the explanation text for each column.
    for ($i = 0; $i -lt $Texts.Length; $i++) {
        $colLetter = [char](65 + $i)
        $ws.Range($colLetter + "2").Value = $Texts[$i]
    }

    # Re-point each dropdown validation so it starts at row 3 (the first
    # real data row) instead of row 2 (now the explanation row). Recreate
    # them in their original order so the validations collection layout is
    # preserved.
    foreach ($v in $Validations) {
        $col = $v.Col
        $formula = $v.Formula

        $oldRange = $ws.Range($col + "2:" + $col + "1000")
        $oldRange.Validation.Delete()

        $newRange = $ws.Range($col + "3:" + $col + "1000")
        $newRange.Validation.Add(3, 1, 1, $formula)
        $newRange.Validation.ShowInput = $false
        $newRange.Validation.ShowError = $false
    }
}

# ---------------------------------------------------------------------------
# CbCR Notifications (columns A-N)
# ---------------------------------------------------------------------------
$cbcrTexts = @(
    "Belgium, France, UK, Germany, etc.",
    "Always / Conditional / Notification Only / Never Required / N/A",
    "Annual / One-Time / Upon Change",
    "UPE / Local CE / One CE for All / Other",
    "Yes / No / Not Specified",
    "Yes / No",
    "Separate Form / Portal / Within CIT Return / BZSt Portal",
    "Form 275.CBC.NOT / DAS2-CbCR / SA / BZSt CbCR Notification",
    "None / CIT Date / FYE-Based / Fixed / Upon Request / With Tax Return",
    "Details (e.g., By 31 Dec following FY, Within 3 months of change)",
    "Valid for FY / Until entity or UPE info changes / Annual / Until change in filing entity",
    "MF / LF / CbCR / Standalone",
    "FY2016, 2017, etc.",
    "Context about notification rules and requirements"
)

$cbcrValidations = @(
    @{ Col = "B"; Formula = '"Always,Conditional,Notification Only,Never Required,N/A"' },
    @{ Col = "C"; Formula = '"Annual,One-Time,Upon Change"' },
    @{ Col = "D"; Formula = '"UPE,Local CE,One CE for All,Other"' },
    @{ Col = "E"; Formula = '"Yes,No,Not Specified"' },
    @{ Col = "F"; Formula = '"Yes,No"' },
    @{ Col = "I"; Formula = '"None,CIT Date,FYE-Based,Fixed,Upon Request,With Tax Return"' },
    @{ Col = "L"; Formula = '"MF,LF,CbCR,Standalone"' }
)

Add-ExplanationRow "CbCR Notifications" "N" $cbcrTexts $cbcrValidations

# ---------------------------------------------------------------------------
# TP Forms (columns A-M)
# ---------------------------------------------------------------------------
$tpFormsTexts = @(
    "Belgium, Spain, Italy, Germany, etc.",
    "Form 275.MF / Form 232 / RS 106 / Transaction Matrix",
    "TP Disclosure / TP Return / MF Summary / LF Summary / CbCR Notification / Other",
    "Always / If MF Required / If LF Required / If MF or LF Required / If CbCR Required / Other",
    "MF / LF / CbCR / Standalone",
    "Summary form with key data / Annual TP informative return / TP disclosure / Structured overview of RPTs",
    "None / CIT Date / FYE-Based / Fixed / Upon Request / With Tax Return",
    "Details (e.g., 31 Dec following FY, Approx 25 Aug, Within 30 days of audit notice)",
    "Days to submit if upon request (30, 14, 10, etc.)",
    "Yes / No",
    "Yes / No / Electronic Timestamp",
    "FY2016, 2010, 2024, etc.",
    "Context about form requirements and special rules"
)

$tpFormsValidations = @(
    @{ Col = "C"; Formula = '"TP Disclosure,TP Return,MF Summary,LF Summary,CbCR Notification,Other"' },
    @{ Col = "D"; Formula = '"Always,If MF Required,If LF Required,If MF or LF Required,If CbCR Required,Other"' },
    @{ Col = "E"; Formula = '"MF,LF,CbCR,Standalone"' },
    @{ Col = "G"; Formula = '"None,CIT Date,FYE-Based,Fixed,Upon Request,With Tax Return"' },
    @{ Col = "J"; Formula = '"Yes,No"' },
    @{ Col = "K"; Formula = '"Yes,No,Electronic Timestamp"' }
)

Add-ExplanationRow "TP Forms" "M" $tpFormsTexts $tpFormsValidations
